# Update loading_percent results for the 380 kV case (Case_1_207, res_line/loading_percent.xlsx)
# New load-flow results replace the previous values in rows 2-25 for columns B,D,E,F,G,H,I,L
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 18.27072548466544
$ws.Range("D2").Value = 3.281823099207152
$ws.Range("E2").Value = 23.64955426240219
$ws.Range("F2").Value = 23.99422088531609
$ws.Range("G2").Value = 30.55378444541328
$ws.Range("H2").Value = 13.22139441892493
$ws.Range("I2").Value = 22.94321346680962
$ws.Range("L2").Value = 11.34947255283348
$ws.Range("B3").Value = 17.64103966924031
$ws.Range("D3").Value = 3.322825915193549
$ws.Range("E3").Value = 22.85656797215459
$ws.Range("F3").Value = 23.57963557676031
$ws.Range("G3").Value = 29.60489905944011
$ws.Range("H3").Value = 13.19082029628856
$ws.Range("I3").Value = 23.08291819917119
$ws.Range("L3").Value = 11.02035033082415
$ws.Range("B4").Value = 17.24286386186172
$ws.Range("D4").Value = 3.349116637291488
$ws.Range("E4").Value = 22.35468758350239
$ws.Range("F4").Value = 23.33477672204104
$ws.Range("G4").Value = 29.02673973814718
$ws.Range("H4").Value = 13.17767010015258
$ws.Range("I4").Value = 23.17813381859774
$ws.Range("L4").Value = 10.81409610229775
$ws.Range("B5").Value = 17.07791033728303
$ws.Range("D5").Value = 3.360111623789527
$ws.Range("E5").Value = 22.14665700845435
$ws.Range("F5").Value = 23.23758845481536
$ws.Range("G5").Value = 28.79274946250067
$ws.Range("H5").Value = 13.17372211771504
$ws.Range("I5").Value = 23.21929493235008
$ws.Range("L5").Value = 10.72911222559358
$ws.Range("B6").Value = 17.05036377656599
$ws.Range("D6").Value = 3.361954351640056
$ws.Range("E6").Value = 22.11190969786611
$ws.Range("F6").Value = 23.22161144881406
$ws.Range("G6").Value = 28.75400777713368
$ws.Range("H6").Value = 13.17315164078093
$ws.Range("I6").Value = 23.22627186318965
$ws.Range("L6").Value = 10.71494786165434
$ws.Range("B7").Value = 17.24064985516499
$ws.Range("D7").Value = 3.349263779292908
$ws.Range("E7").Value = 22.35189586440757
$ws.Range("F7").Value = 23.33345530629209
$ws.Range("G7").Value = 29.02357686372091
$ws.Range("H7").Value = 13.17761114953167
$ws.Range("I7").Value = 23.17867939350804
$ws.Range("L7").Value = 10.81295360172462
$ws.Range("B8").Value = 18.05612566906511
$ws.Range("D8").Value = 3.295730072830162
$ws.Range("E8").Value = 23.37938921145286
$ws.Range("F8").Value = 23.84935421434361
$ws.Range("G8").Value = 30.22597309147599
$ws.Range("H8").Value = 13.20968436765264
$ws.Range("I8").Value = 22.98941742372466
$ws.Range("L8").Value = 11.23692381713158
$ws.Range("B9").Value = 19.55548028252181
$ws.Range("D9").Value = 3.199551671027258
$ws.Range("E9").Value = 25.26535270026964
$ws.Range("F9").Value = 24.93085108439016
$ws.Range("G9").Value = 32.59864751351009
$ws.Range("H9").Value = 13.31722462430591
$ws.Range("I9").Value = 22.6937220161969
$ws.Range("L9").Value = 12.03057390096252
$ws.Range("B10").Value = 20.58665137297375
$ws.Range("D10").Value = 3.134192014184013
$ws.Range("E10").Value = 26.56055978296185
$ws.Range("F10").Value = 25.75818020722365
$ws.Range("G10").Value = 34.32512755823367
$ws.Range("H10").Value = 13.42335377238371
$ws.Range("I10").Value = 22.52324856689592
$ws.Range("L10").Value = 12.58501909313759
$ws.Range("B11").Value = 21.03881364591554
$ws.Range("D11").Value = 3.105596484942549
$ws.Range("E11").Value = 27.12810099249745
$ws.Range("F11").Value = 26.13956902090479
$ws.Range("G11").Value = 35.10208941304388
$ws.Range("H11").Value = 13.47745776280013
$ws.Range("I11").Value = 22.45601632306492
$ws.Range("L11").Value = 12.83000676778023
$ws.Range("B12").Value = 21.20748616592185
$ws.Range("D12").Value = 3.094930637125595
$ws.Range("E12").Value = 27.33975590443806
$ws.Range("F12").Value = 26.28454316102055
$ws.Range("G12").Value = 35.39473709564388
$ws.Range("H12").Value = 13.49877439463839
$ws.Range("I12").Value = 22.43205396641643
$ws.Range("L12").Value = 12.92166286170861
$ws.Range("B13").Value = 21.17127472195069
$ws.Range("D13").Value = 3.097220496894472
$ws.Range("E13").Value = 27.29431939179038
$ws.Range("F13").Value = 26.25329861085481
$ws.Range("G13").Value = 35.33178541493637
$ws.Range("H13").Value = 13.49414678867133
$ws.Range("I13").Value = 22.43714789625071
$ws.Range("L13").Value = 12.90197380978595
$ws.Range("B14").Value = 21.05274219511123
$ws.Range("D14").Value = 3.104715743955649
$ws.Range("E14").Value = 27.14558009268146
$ws.Range("F14").Value = 26.15148607375411
$ws.Range("G14").Value = 35.12619913727855
$ws.Range("H14").Value = 13.47919494704567
$ws.Range("I14").Value = 22.45401484664599
$ws.Range("L14").Value = 12.83757012914245
$ws.Range("B15").Value = 20.97980206835605
$ws.Range("D15").Value = 3.109327958015228
$ws.Range("E15").Value = 27.0540443062685
$ws.Range("F15").Value = 26.08918958151166
$ws.Range("G15").Value = 35.00005693023686
$ws.Range("H15").Value = 13.47014412075194
$ws.Range("I15").Value = 22.4645416937023
$ws.Range("L15").Value = 12.79797364755068
$ws.Range("B16").Value = 20.55675054007853
$ws.Range("D16").Value = 3.136083605439961
$ws.Range("E16").Value = 26.52302091762119
$ws.Range("F16").Value = 25.73334191267009
$ws.Range("G16").Value = 34.27415026083295
$ws.Range("H16").Value = 13.41993444245111
$ws.Range("I16").Value = 22.52785111359536
$ws.Range("L16").Value = 12.56885641527798
$ws.Range("B17").Value = 20.29279898243952
$ws.Range("D17").Value = 3.152787920967959
$ws.Range("E17").Value = 26.19159864126607
$ws.Range("F17").Value = 25.51620527872598
$ws.Range("G17").Value = 33.82640604190147
$ws.Range("H17").Value = 13.39061864311465
$ws.Range("I17").Value = 22.56934148562263
$ws.Range("L17").Value = 12.42639190997326
$ws.Range("B18").Value = 20.13939423741886
$ws.Range("D18").Value = 3.162502860465736
$ws.Range("E18").Value = 25.99894258133721
$ws.Range("F18").Value = 25.39179786536163
$ws.Range("G18").Value = 33.56810198444285
$ws.Range("H18").Value = 13.37430595019219
$ws.Range("I18").Value = 22.59417594441225
$ws.Range("L18").Value = 12.34377336481638
$ws.Range("B19").Value = 20.08718534060164
$ws.Range("D19").Value = 3.165810585129351
$ws.Range("E19").Value = 25.93336854227713
$ws.Range("F19").Value = 25.34976410107009
$ws.Range("G19").Value = 33.4805233515813
$ws.Range("H19").Value = 13.36887728261521
$ws.Range("I19").Value = 22.60275069950539
$ws.Range("L19").Value = 12.31568634250699
$ws.Range("B20").Value = 20.32106220623571
$ws.Range("D20").Value = 3.150998643939847
$ws.Range("E20").Value = 26.22709040219952
$ws.Range("F20").Value = 25.53927103808615
$ws.Range("G20").Value = 33.87415193900101
$ws.Range("H20").Value = 13.39368260181929
$ws.Range("I20").Value = 22.56482424453501
$ws.Range("L20").Value = 12.44162810984053
$ws.Range("B21").Value = 21.08762818592665
$ws.Range("D21").Value = 3.102509800164118
$ws.Range("E21").Value = 27.18935800141838
$ws.Range("F21").Value = 26.18137730123632
$ws.Range("G21").Value = 35.18663014322978
$ws.Range("H21").Value = 13.48356425881114
$ws.Range("I21").Value = 22.44901987945291
$ws.Range("L21").Value = 12.85651788643754
$ws.Range("B22").Value = 21.57370665654364
$ws.Range("D22").Value = 3.071767270678726
$ws.Range("E22").Value = 27.79919456498944
$ws.Range("F22").Value = 26.6041627280507
$ws.Range("G22").Value = 36.03511103958374
$ws.Range("H22").Value = 13.54713137881086
$ws.Range("I22").Value = 22.38206693296046
$ws.Range("L22").Value = 13.12114117980481
$ws.Range("B23").Value = 21.31567751533461
$ws.Range("D23").Value = 3.088088686149137
$ws.Range("E23").Value = 27.47550121180433
$ws.Range("F23").Value = 26.37828356324545
$ws.Range("G23").Value = 35.58322182191544
$ws.Range("H23").Value = 13.51276649511654
$ws.Range("I23").Value = 22.41699763830102
$ws.Range("L23").Value = 12.98052729483392
$ws.Range("B24").Value = 20.30828955008722
$ws.Range("D24").Value = 3.151807229016943
$ws.Range("E24").Value = 26.21105115967557
$ws.Range("F24").Value = 25.52884166404684
$ws.Range("G24").Value = 33.85256876322151
$ws.Range("H24").Value = 13.39229569924053
$ws.Range("I24").Value = 22.56686343437395
$ws.Range("L24").Value = 12.43474204141219
$ws.Range("B25").Value = 19.16160934735012
$ws.Range("D25").Value = 3.224634820242035
$ws.Range("E25").Value = 24.77026102686983
$ws.Range("F25").Value = 24.63180939153224
$ws.Range("G25").Value = 31.95796561164866
$ws.Range("H25").Value = 13.28335414787528
$ws.Range("I25").Value = 22.76555731188869
$ws.Range("L25").Value = 11.82052693028794
